$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ForeignAccountData")

# ---------------------------------------------------------------------
# 1) Numeric / text value updates (POR account number corrections)
# ---------------------------------------------------------------------

# Rows where W, AC and AI held the old account number 926091891 (as a
# plain number) and must become 28048100718.
$rowsPlainAcct = @(11,12,14,15,16,17)
foreach ($r in $rowsPlainAcct) {
    $ws.Cells.Item($r, 23).Value  = 28048100718   # W
    $ws.Cells.Item($r, 29).Value  = 28048100718   # AC
    $ws.Cells.Item($r, 35).Value  = 28048100718   # AI
}

# Row 5 is special: AC5 is formatted as Text (numFmt 49) and the new
# account number is entered there as text, not as a number.
$ws.Cells.Item(5, 23).Value = 28048100718            # W5
$ws.Cells.Item(5, 35).Value = 28048100718            # AI5
$ws.Range("AC5").NumberFormat = "@"
$ws.Range("AC5").Value = "28048100718"

# AY column: 32400033618 -> 42010250514 on the same set of rows as W/AI
# (rows 5, 11, 12, 14-17)
$rowsAY = @(5,11,12,14,15,16,17)
foreach ($r in $rowsAY) {
    $ws.Cells.Item($r, 51).Value = 42010250514   # AY
}

# Row 6 AY changes to a different, unrelated POR number.
$ws.Cells.Item(6, 51).Value = 42010039950

# ---------------------------------------------------------------------
# 2) Row 6 formatting: the whole row was given a white "Background 1"
#    solid fill (customFormat row), as when selecting the row and
#    applying a fill color from the palette.
# ---------------------------------------------------------------------

$row6Cols = @(1,2,3,4,6,7,8,9,10,11,12,13,14,16,17,18,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,42,43,44,45,49,51,52,53,69,70,72,74,75,78,80,81,84,85,86,89,91,92,94,99,147,148,149)

$firstCell = $ws.Cells.Item(6, $row6Cols[0])
$firstCell.Interior.ThemeColor = 2
$firstCell.Interior.TintAndShade = 0
$firstCell.Copy()

foreach ($c in $row6Cols) {
    if ($c -ne $row6Cols[0]) {
        $ws.Cells.Item(6, $c).PasteSpecial(-4122)
    }
}

# AC6 (column 29) keeps its existing Text number format (numFmtId 49)
# in addition to picking up the new fill.
$ws.Cells.Item(6, 29).NumberFormat = "@"

$excel.CutCopyMode = 0
